$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (Worksheets.Item(2))
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

# Status column -> "Handed back: in sync with en-US"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("B3").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime column (was the "0001-01-01 00:00:00" placeholder)
$ws.Range("G2").Value = "2016-01-26 12:18:00"
$ws.Range("G3").Value = "2016-01-26 12:18:00"

# Rebuild all hyperlinks (existing ones + the two new "Latest Target File" /
# "Latest Handback File" columns) so they line up in the expected order.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ca863c34c06a2da33902fd37cd6f21c3253cc0a1/e2e/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a9e658be81d8a71e56562cc4a5897933818b65e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/ca863c34c06a2da33902fd37cd6f21c3253cc0a1/e2e/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a9e658be81d8a71e56562cc4a5897933818b65e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ca863c34c06a2da33902fd37cd6f21c3253cc0a1/e2e/ffff4f896d8c-d0bd-4be3-aacb-602751cc319c.md", "", "", "ffff4f896d8c-d0bd-4be3-aacb-602751cc319c.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a9e658be81d8a71e56562cc4a5897933818b65e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/ca863c34c06a2da33902fd37cd6f21c3253cc0a1/e2e/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a9e658be81d8a71e56562cc4a5897933818b65e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ca863c34c06a2da33902fd37cd6f21c3253cc0a1/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de" (Worksheets.Item(3))
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(3)

# Status column -> "Handed back: in sync with en-US"
$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("B3").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime column (was the "0001-01-01 00:00:00" placeholder)
$ws2.Range("G2").Value = "2016-01-26 12:18:21"
$ws2.Range("G3").Value = "2016-01-26 12:18:21"

# Rebuild all hyperlinks (existing ones + the two new "Latest Target File" /
# "Latest Handback File" columns) so they line up in the expected order.
$ws2.Hyperlinks.Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ca863c34c06a2da33902fd37cd6f21c3253cc0a1/e2e/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70fbc5520a1830457cc4ae56152d86c9e920b2d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/ca863c34c06a2da33902fd37cd6f21c3253cc0a1/e2e/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70fbc5520a1830457cc4ae56152d86c9e920b2d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ca863c34c06a2da33902fd37cd6f21c3253cc0a1/e2e/ffff4f896d8c-d0bd-4be3-aacb-602751cc319c.md", "", "", "ffff4f896d8c-d0bd-4be3-aacb-602751cc319c.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70fbc5520a1830457cc4ae56152d86c9e920b2d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/ca863c34c06a2da33902fd37cd6f21c3253cc0a1/e2e/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70fbc5520a1830457cc4ae56152d86c9e920b2d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ca863c34c06a2da33902fd37cd6f21c3253cc0a1/.localization-config", "", "", ".localization-config")

Write-Output "Generate Report for handback - done"
